# Fix booking-system workbook: align the row-3 booking with the same
# person as row 2 (Huzaifa RAGHAV), lower-case the AM/PM time strings,
# and add a "Reminder" column (H) so that your_bookings / cancel e-mail
# lookups (which depend on a spaceNameWithSpaces-style match) stop
# breaking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 booking belongs to the same person as row 2 (Huzaifa RAGHAV).
$ws.Range("A3").Value = "Huzaifa RAGHAV"

# Update B3's text (leading "'" keeps the existing quote-prefixed Hyperlink
# cell style instead of minting a fresh, near-duplicate style entry).
$ws.Range("B3").Value = "'ragha59105@gapps.uwcsea.edu.sg"

# Normalise the AM/PM time strings to lowercase.
$ws.Range("D2").Value = "12:45pm"
$ws.Range("E2").Value = "1:15pm"
$ws.Range("D3").Value = "12:45pm"
$ws.Range("E3").Value = "1:15pm"

# New "Reminder" column with the reminder lead-time per booking.
$ws.Range("H1").Value = "Reminder "
$ws.Range("H2").Value = "12h"
$ws.Range("H3").Value = "4h"

# Match the cursor position recorded in the edited workbook.
$ws.Range("E9").Select()
